# Update the "dSF" (column F) values on Sheet1 with re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -1
    3  = 5
    4  = 1
    5  = -1
    6  = 2
    7  = -2
    9  = 5
    10 = 7
    12 = 2
    13 = 4
    15 = 5
    16 = 5
    17 = -1
    18 = 4
    19 = 0
    20 = -1
    22 = 9
    25 = -3
    26 = 1
    27 = 4
    28 = 4
    29 = 1
    31 = -1
    33 = 1
    34 = -2
    35 = -1
    36 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
